$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) " and its effects on information cleavages has received significant
#    scholarly interest. This work looks at incidental exposure and social
#    inequality " -> " has received significant scholarly interest.
#    Scholarship has examined serendipitous news exposure and information
#    inequality "
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    " and its effects on information cleavages has received significant scholarly interest. This work looks at incidental exposure and social inequality ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " has received significant scholarly interest. Scholarship has examined serendipitous news exposure and information inequality ",
    2) | Out-Null

# ------------------------------------------------------------------
# 2) "exposure to politically relevant information" -> "access to politically relevant information"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "in the form of exposure to politically relevant information (",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in the form of access to politically relevant information (",
    2) | Out-Null

# ------------------------------------------------------------------
# 3) "Two frames dominate thinking in this area: compensatory effects of
#    information heterogeneity to engage an otherwise disinterested public
#    (), or stratification via the 'Matthew " -> new text
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Two frames dominate thinking in this area: compensatory effects of information heterogeneity to engage an otherwise disinterested public (), or stratification via the " + [char]0x2018 + "Matthew ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The two dominant frameworks for thinking in this area are based on " + [char]0x2018 + "compensatory" + [char]0x2019 + " effects, or the ability of information heterogeneity to engage an otherwise disinterested public (Ahmadi & Wohn, 2018), or stratification effects via the " + [char]0x2018 + "Matthew ",
    2) | Out-Null

# ------------------------------------------------------------------
# 4) "...as the 'tune out' of..." -> "...as they 'tune out' of..."
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "social gratifications as the " + [char]0x2018 + "tune out" + [char]0x2019 + " of political life",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "social gratifications as they " + [char]0x2018 + "tune out" + [char]0x2019 + " of political life",
    2) | Out-Null

# ------------------------------------------------------------------
# 5) "Empirical studies do offer a clear picture, as evidence of equalization
#    and stratification are similarly possible, depending on various
#    socio-technical conditions. " -> long rewritten/expanded paragraph
# ------------------------------------------------------------------
$old5 = "Empirical studies do offer a clear picture, as evidence of equalization and stratification are similarly possible, depending on various socio-technical conditions. "
$new5 = "Empirical findings offer a complicated picture. Studies provide evidence that equalization and stratification are similarly possible, depending on various socio-technical conditions. First, Fletcher and Nielsen (2018) show strong and convincing evidence for equalization effects in terms of news exposure. Using survey data from four Western countries (Italy, Australia, United Kingdom, United States), they find that people who use social media for purposes other than news are exposed to significantly more online news sources, and the effect is stronger among those with lower levels of political interest. In another cross-national sample, semi-structured interviews reveal that this effect is due, at least in part, to stumbling across topics of potential interest based on activities of others on the platform, thus pulling the otherwise disengaged into an information/engagement feedback loop (Mitchelstein et al., 2020). These findings qualify the nature of incidentality, as dependent upon not only information heterogeneity" + [char]0x2014 + "as traditionally theorized (e.g., Tewksbury et al., 2001)" + [char]0x2014 + "but also network size and diversity, as larger networks increase the chances of both incidental and purposeful news use (Barnidge, 2021)."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

# ------------------------------------------------------------------
# 6) Structural changes after the "Empirical ..." paragraph:
#    - the following (currently empty) paragraph gets "Second, " text
#    - the paragraph after that gets spacing(480,auto)/firstLine(720) added
#    - the empty paragraph immediately before "II. Stratificational effects"
#      and that heading paragraph itself are both removed
# ------------------------------------------------------------------
$empiricalIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Empirical findings offer")) {
        $empiricalIdx = $i
        break
    }
}
Write-Output "empiricalIdx=$empiricalIdx"

$secondPara = $d.Paragraphs.Item($empiricalIdx + 1)
$secondPara.Range.Text = "Second, "
$secondPara.Range.Font.Name = "Times New Roman"
$secondPara.Range.Font.NameAscii = "Times New Roman"
$secondPara.Range.Font.NameBi = "Times New Roman"
$secondPara.Range.Font.NameOther = "Times New Roman"

$nextPara = $d.Paragraphs.Item($empiricalIdx + 2)
$nextPara.Format.LineSpacingRule = 2
$nextPara.Format.FirstLineIndent = 36

# locate the "II. Stratificational effects" paragraph and the empty one before it
$stratIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("II. Stratificational")) {
        $stratIdx = $i
        break
    }
}
Write-Output "stratIdx=$stratIdx"

$delStart = $d.Paragraphs.Item($stratIdx - 1).Range.Start
$delEnd = $d.Paragraphs.Item($stratIdx).Range.End
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()

Write-Output "Stage 1 done"
